$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.061.45'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '2.259.50'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.73'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.525'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.04%  '
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.86'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.42%  '
$ws.Range('E11').Value = '  -2.09%  '
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.82'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.73%  '
$ws.Range('D14').Value = '2.610.05'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.59'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '2.261.25'
$ws.Range('E16').Value = '  -1.44%  '
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('D18').Value = '41.889.57'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.93%  '
$ws.Range('D20').Value = '0.0₃0899'
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.47'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.69%  '
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.64%  '
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '164.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.21'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.69%  '
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.14'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.66'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.06%  '
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('E38').Value = '  -4.85%  '
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('E41').Value = '  -1.51%  '
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('D43').Value = '1.946.69'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.75%  '
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('D49').Value = '2.483.10'
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.85'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '91.88'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.10%  '
